$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.7025
$ws.Range("C4").Value = -11.90690000000001

$ws.Range("A7").Value = -19.21589999999998

$ws.Range("D10").Value = -8.181099999999997

$ws.Range("C12").Value = -10.37459999999999

$ws.Range("D13").Value = -8.454299999999995

$ws.Range("A16").Value = -22.01390000000001

$ws.Range("C18").Value = -12.27209999999999

$ws.Range("C19").Value = -10.9974

$ws.Range("C20").Value = -12.07699999999999

$ws.Range("A28").Value = -21.912

$ws.Range("A29").Value = -21.08849999999997

$ws.Range("D30").Value = -6.887799999999996

$ws.Range("C31").Value = -13.1193

$ws.Range("A32").Value = -21.11040000000001

$ws.Range("A40").Value = -20.55720000000001
$ws.Range("C40").Value = -11.56910000000002
$ws.Range("D40").Value = -7.545999999999998

$ws.Range("C42").Value = -11.20130000000001

$ws.Range("D44").Value = -6.512800000000004

$ws.Range("C47").Value = -12.11269999999999

$ws.Range("C48").Value = -11.5942

$ws.Range("A52").Value = -22.22409999999999

$ws.Range("A57").Value = -22.2093

$ws.Range("C63").Value = -10.216

$ws.Range("C64").Value = -10.77289999999999

$ws.Range("A66").Value = -21.40329999999999

$ws.Range("C76").Value = -11.9184

$ws.Range("C81").Value = -14.0175

$ws.Range("C89").Value = -13.6642
$ws.Range("D89").Value = -8.631199999999993

$ws.Range("D91").Value = -7.9237

$ws.Range("C94").Value = -10.7814

$ws.Range("A100").Value = -22.03310000000003
